$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row right-answer count
$ws.Range("B11").Value = 5

# Update "Total" row correct marks and correct/total string
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
